# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" table (rows 16-27) is re-sequenced: a new worker
# (CARLOS ANDRES CORREDOR MONTERROSA, CC 1020738585) is interleaved with
# the existing worker (ANA DEL ROSARIO CORREDOR MONTERROSA, CC 1020720479)
# for periods 2110/2111, and ANA's remaining periods are reordered
# chronologically (2201, 2203-2209).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: ANA - period 2110
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1020720479"
$ws.Range("D16").Value = "ANA DEL ROSARIO CORREDOR MONTERROSA"
$ws.Range("E16").Value = "2110"
$ws.Range("F16").Value = 36341
$ws.Range("G16").Value = 1000000

# Row 17: CARLOS - period 2110
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1020738585"
$ws.Range("D17").Value = "CARLOS ANDRES CORREDOR MONTERROSA"
$ws.Range("E17").Value = "2110"
$ws.Range("F17").Value = 36341
$ws.Range("G17").Value = 1000000

# Row 18: ANA - period 2111
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1020720479"
$ws.Range("D18").Value = "ANA DEL ROSARIO CORREDOR MONTERROSA"
$ws.Range("E18").Value = "2111"
$ws.Range("F18").Value = 36341
$ws.Range("G18").Value = 1000000

# Row 19: CARLOS - period 2111
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1020738585"
$ws.Range("D19").Value = "CARLOS ANDRES CORREDOR MONTERROSA"
$ws.Range("E19").Value = "2111"
$ws.Range("F19").Value = 36341
$ws.Range("G19").Value = 1000000

# Row 20: ANA - period 2201
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1020720479"
$ws.Range("D20").Value = "ANA DEL ROSARIO CORREDOR MONTERROSA"
$ws.Range("E20").Value = "2201"
$ws.Range("F20").Value = 36341
$ws.Range("G20").Value = 1000000

# Row 21: ANA - period 2203
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1020720479"
$ws.Range("D21").Value = "ANA DEL ROSARIO CORREDOR MONTERROSA"
$ws.Range("E21").Value = "2203"
$ws.Range("F21").Value = 40000
$ws.Range("G21").Value = 1000000

# Row 22: ANA - period 2204
$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1020720479"
$ws.Range("D22").Value = "ANA DEL ROSARIO CORREDOR MONTERROSA"
$ws.Range("E22").Value = "2204"
$ws.Range("F22").Value = 40000
$ws.Range("G22").Value = 1000000

# Row 23: ANA - period 2205
$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1020720479"
$ws.Range("D23").Value = "ANA DEL ROSARIO CORREDOR MONTERROSA"
$ws.Range("E23").Value = "2205"
$ws.Range("F23").Value = 40000
$ws.Range("G23").Value = 1000000

# Row 24: ANA - period 2206
$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "1020720479"
$ws.Range("D24").Value = "ANA DEL ROSARIO CORREDOR MONTERROSA"
$ws.Range("E24").Value = "2206"
$ws.Range("F24").Value = 40000
$ws.Range("G24").Value = 1000000

# Row 25: ANA - period 2207
$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "1020720479"
$ws.Range("D25").Value = "ANA DEL ROSARIO CORREDOR MONTERROSA"
$ws.Range("E25").Value = "2207"
$ws.Range("F25").Value = 40000
$ws.Range("G25").Value = 1000000

# Row 26: ANA - period 2208
$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "1020720479"
$ws.Range("D26").Value = "ANA DEL ROSARIO CORREDOR MONTERROSA"
$ws.Range("E26").Value = "2208"
$ws.Range("F26").Value = 40000
$ws.Range("G26").Value = 1000000

# Row 27: ANA - period 2209
$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "1020720479"
$ws.Range("D27").Value = "ANA DEL ROSARIO CORREDOR MONTERROSA"
$ws.Range("E27").Value = "2209"
$ws.Range("F27").Value = 34666
$ws.Range("G27").Value = 1000000
